$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture current values before shifting so we don't clobber them
$valC1 = $ws.Range("C1").Value2
$valD1 = $ws.Range("D1").Value2
$valE1 = $ws.Range("E1").Value2

$valD2 = $ws.Range("D2").Value2

# Header row: shift max/prediction/rejection-f left by one column (C<-D, D<-E, E<-C)
$ws.Range("C1").Value = $valD1
$ws.Range("D1").Value = $valE1
$ws.Range("E1").Value = $valC1

# Data row 2: C becomes the species text (same as D), D stays same species text, E becomes 1
$ws.Range("C2").Value = $valD2
$ws.Range("D2").Value = $valD2
$ws.Range("E2").Value = 1
